$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 62

$valuesA = @()
$valuesB = @()

for ($r = 1; $r -le $lastRow; $r++) {
    $valuesA += $ws.Cells.Item($r, 1).Value2
    $valuesB += $ws.Cells.Item($r, 2).Value2
}

for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $valuesB[$r - 1]
    $ws.Cells.Item($r, 2).Value2 = $valuesA[$r - 1]
}

# Row 59 also carries a distinct highlight font color that travels with the
# value it was attached to (it was on column A, now belongs to column B).
$colorA59 = $ws.Cells.Item(59, 1).Font.Color
$colorB59 = $ws.Cells.Item(59, 2).Font.Color
$ws.Cells.Item(59, 1).Font.Color = $colorB59
$ws.Cells.Item(59, 2).Font.Color = $colorA59
